# Auto-generated Excel COM-interop script to update cryptos worksheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "25.842.25"
$ws.Range("E2").Value = "  +0.18%  "

Set-TextCell "D3" "1.630.42"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("E4").Value = "  +0.62%  "

Set-TextCell "D5" "214.13"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("E7").Value = "  +0.50%  "

$ws.Range("E8").Value = "  -0.40%  "

Set-TextCell "D9" "0.0632"
$ws.Range("E9").Value = "  +0.20%  "

Set-TextCell "D10" "19.53"
$ws.Range("E10").Value = "  -0.46%  "

Set-TextCell "D11" "0.0790"
$ws.Range("E11").Value = "  +0.34%  "

Set-TextCell "D12" "1.856.11"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("E13").Value = "  -0.07%  "

Set-TextCell "D14" "1.632.91"
$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("E15").Value = "  -1.32%  "

Set-TextCell "D16" "0.0₃0754"
$ws.Range("E16").Value = "  -0.47%  "

Set-TextCell "D17" "62.59"
$ws.Range("E17").Value = "  -0.10%  "

Set-TextCell "D18" "25.848.29"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D20" "193.11"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D21" "4.38"
$ws.Range("E21").Value = "  -1.05%  "

Set-TextCell "D22" "9.91"
$ws.Range("E22").Value = "  -0.06%  "

Set-TextCell "D23" "6.25"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("E24").Value = "  +1.16%  "

Set-TextCell "D25" "143.18"
$ws.Range("E25").Value = "  +0.62%  "

$ws.Range("E26").Value = "  +0.52%  "

$ws.Range("E27").Value = "  +2.66%  "

Set-TextCell "D28" "6.83"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("E34").Value = "  -2.10%  "

$ws.Range("E35").Value = "  +1.51%  "

$ws.Range("E36").Value = "  -0.27%  "

Set-TextCell "D37" "1.136.63"
$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("E39").Value = "  -1.10%  "

$ws.Range("E40").Value = "  +0.48%  "

$ws.Range("E41").Value = "  +0.50%  "

Set-TextCell "D42" "99.08"
$ws.Range("E42").Value = "  -1.32%  "

Set-TextCell "D43" "5.43"
$ws.Range("E43").Value = "  -3.03%  "

Set-TextCell "D44" "0.794"
$ws.Range("E44").Value = "  -0.54%  "

Set-TextCell "D45" "1.765.64"

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D46" "56.19"
$ws.Range("E46").Value = "  +1.55%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D47" "0.0527"
$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D48" "1.44"
$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D49" "0.415"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D50" "7.61"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D51" "0.0957"
$ws.Range("E51").Value = "  +0.35%  "
